# The sheet originally held a pandas DataFrame written with its integer
# row-index exported as column A (bold/boxed header style carried over onto
# the index cells), and the real "a/b/c/d" data shifted one column to the
# right (B:E). The re-export drops the index column, so everything just
# shifts left by one column: B:E -> A:D.
#
# Deleting column A reproduces exactly that: the old header row
# (B1:E1 = "a","b","c","d", bold/boxed style) becomes A1:D1, and the old
# data rows (B2:E5) become A2:D5, while the previous index column (A, which
# carried the header style) is removed entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:A").Delete()
